$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new Mac-Address rows appended to the reg_center_user_machine_h master data
$newRows = @(
    @(10001, 110030, 10030, "eng", $true, "superadmin", "now()", "now()"),
    @(10001, 110031, 10031, "eng", $true, "superadmin", "now()", "now()")
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowData[$j]
    }
}

# Move the view/selection the same way the author left the sheet: scrolled down to
# the new rows, with F30 as the active cell.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F30").Select()
